# "Fixed total used capacity and added all of the variables descriptions"
#
# 1. company sheet: add a new variable row ("average delay time in
#    dropping out") and widen column A to fit the new, longer label.
# 2. demand sheet: remove the stray "goal" row (it belonged logically with
#    the other cost-related goal, not demand).
# 3. cost sheet: fix the available capacity / total used capacity values
#    (hubs location variable, routing & consolidation efficiency, total
#    operation cost of one hub) and make this sheet the active tab.

$wb = $excel.ActiveWorkbook

# --- company sheet ---------------------------------------------------
$ws1 = $wb.Worksheets.Item("company")
$ws1.Range("A5").Value = "average delay time in dropping out"
$ws1.Range("B5").Value = 36
$ws1.Range("C5").Value = "month"
$ws1.Columns.Item(1).ColumnWidth = 29.75
$ws1.Range("A3").Select() | Out-Null

# --- demand sheet ------------------------------------------------------
$ws3 = $wb.Worksheets.Item("demand")
$ws3.Rows.Item(3).Delete() | Out-Null
$ws3.Range("B9").Select() | Out-Null

# --- cost sheet ----------------------------------------------------------
$ws4 = $wb.Worksheets.Item("cost")
$ws4.Range("B1").Value = 0.1
$ws4.Range("B2").Value = 0.7
$ws4.Range("B3").Value = 100
# Activate last so it becomes the workbook's active tab / selected sheet.
$ws4.Activate() | Out-Null
$ws4.Range("B2").Select() | Out-Null
